$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# ---------------------------------------------------------------------------
# Build the two new "header" cell styles once, on scratch cells far outside
# the used range, then propagate them with a format-only paste. (Building a
# border style incrementally with several sequential Borders(...).LineStyle
# assignments on many different cells in a row can make the engine mint
# extra, unused intermediate cell styles; a single isolated build followed by
# PasteSpecial(xlPasteFormats) avoids that and keeps the style table minimal,
# matching styles actually used: top+bottom only, and top+right+bottom.)
# ---------------------------------------------------------------------------
$scratchTopBottom = $ws1.Range("Z1")
$scratchTopRightBottom = $ws1.Range("Z2")

$scratchTopBottom.ClearFormats()
$scratchTopBottom.Borders.Item(8).LineStyle = 1
$scratchTopBottom.Borders.Item(9).LineStyle = 1

$scratchTopRightBottom.ClearFormats()
$scratchTopRightBottom.Borders.Item(8).LineStyle = 1
$scratchTopRightBottom.Borders.Item(10).LineStyle = 1
$scratchTopRightBottom.Borders.Item(9).LineStyle = 1

# quality_comparison: C1 / D1
$scratchTopBottom.Copy()
$ws1.Range("C1").PasteSpecial(-4122)
$scratchTopRightBottom.Copy()
$ws1.Range("D1").PasteSpecial(-4122)

# computational_comparison: C1 / D1 and F1 / G1
$scratchTopBottom.Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)
$scratchTopRightBottom.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

# Clean up the scratch cells and the clipboard marching ants.
$ws1.Range("Z1:Z2").ClearContents()
$ws1.Range("Z1:Z2").ClearFormats()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Anonymize "fedcore" -> "approach"
# ---------------------------------------------------------------------------
$ws1.Range("C2").Value2 = "approach"
$ws2.Range("C2").Value2 = "approach"
$ws2.Range("F2").Value2 = "approach"

# ---------------------------------------------------------------------------
# Remove the stray empty cell G5 on computational_comparison
# ---------------------------------------------------------------------------
$ws2.Range("G5").ClearContents()

Write-Host "edits applied"
